$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1433
$ws.Range("J17").Value = 1433
$ws.Range("L17").Value = 4299
$ws.Range("N17").Value = -4635
$ws.Range("H28").Value = 2329
$ws.Range("I28").Value = 2329
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 2329
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = ""
$ws.Range("N28").Value = -1844
$ws.Range("H40").Value = 2989.6
$ws.Range("J40").Value = 2989.6
$ws.Range("L40").Value = 2989.6
$ws.Range("N40").Value = -3339.6
$ws.Range("H62").Value = 4927.5
$ws.Range("I62").Value = 4141.6665
$ws.Range("K62").Value = 4141.6665
$ws.Range("M62").Value = -3517.6665
$ws.Range("H65").Value = 4927.5
$ws.Range("I65").Value = 4141.6665
$ws.Range("K65").Value = 20708.3325
$ws.Range("M65").Value = -17588.3325
$ws.Range("H69").Value = 9156.25
$ws.Range("I69").Value = 8950
$ws.Range("J69").Value = 9500
$ws.Range("K69").Value = 26850
$ws.Range("L69").Value = 28500
$ws.Range("M69").Value = -25976
$ws.Range("N69").Value = -30248
$ws.Range("H72").Value = 9156.25
$ws.Range("I72").Value = 8950
$ws.Range("J72").Value = 9500
$ws.Range("K72").Value = 80550
$ws.Range("L72").Value = 85500
$ws.Range("M72").Value = -76182
$ws.Range("N72").Value = -94236
$ws.Range("H98").Value = 30294.1
$ws.Range("I98").Value = 47516.582
$ws.Range("J98").Value = 4460.375
$ws.Range("K98").Value = 47516.582
$ws.Range("L98").Value = 4460.375
$ws.Range("M98").Value = -46018.582
$ws.Range("N98").Value = -7456.375
$ws.Range("H106").Value = 6178424
$ws.Range("I106").Value = 8235833
$ws.Range("J106").Value = 6197
$ws.Range("K106").Value = 8235833
$ws.Range("L106").Value = 6197
$ws.Range("M106").Value = -8235202
$ws.Range("N106").Value = -7459
$ws.Range("H112").Value = 2423.1052
$ws.Range("J112").Value = 2085.9333
$ws.Range("L112").Value = 6257.7999
$ws.Range("N112").Value = -8473.7999
$ws.Range("H121").Value = 2753.8333
$ws.Range("J121").Value = 2753.8333
$ws.Range("L121").Value = 8261.499899999999
$ws.Range("N121").Value = -11755.4999
$ws.Range("H122").Value = 30294.1
$ws.Range("I122").Value = 47516.582
$ws.Range("J122").Value = 4460.375
$ws.Range("K122").Value = 142549.746
$ws.Range("L122").Value = 13381.125
$ws.Range("M122").Value = -140099.746
$ws.Range("N122").Value = -18281.125
$ws.Range("H132").Value = 3336920.2
$ws.Range("I132").Value = 3601.111
$ws.Range("K132").Value = 10803.333
$ws.Range("M132").Value = -8273.332999999999
$ws.Range("H138").Value = 294863.06
$ws.Range("I138").Value = 659898.75
$ws.Range("K138").Value = 1979696.25
$ws.Range("M138").Value = -1974556.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 412800
$ws.Range("I45").Value = 1015000
$ws.Range("K45").Value = 1015000
$ws.Range("M45").Value = -1014623
$ws.Range("H122").Value = 1118080.1
$ws.Range("I122").Value = 6385.6875
$ws.Range("K122").Value = 19157.0625
$ws.Range("M122").Value = -16707.0625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4658.4614
$ws.Range("I20").Value = 2845.375
$ws.Range("K20").Value = 2845.375
$ws.Range("M20").Value = -2598.375
$ws.Range("H86").Value = 5426.9565
$ws.Range("I86").Value = 6423.125
$ws.Range("K86").Value = 6423.125
$ws.Range("M86").Value = -5300.125
$ws.Range("H89").Value = 5426.9565
$ws.Range("I89").Value = 6423.125
$ws.Range("K89").Value = 32115.625
$ws.Range("M89").Value = -26499.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 159428
$ws.Range("I99").Value = 296272.88
$ws.Range("K99").Value = 296272.88
$ws.Range("M99").Value = -294774.88
$ws.Range("H122").Value = 7224.579
$ws.Range("I122").Value = 7224.579
$ws.Range("K122").Value = 21673.737
$ws.Range("M122").Value = -19223.737
$ws.Range("H126").Value = 159428
$ws.Range("I126").Value = 296272.88
$ws.Range("K126").Value = 888818.64
$ws.Range("M126").Value = -886348.64

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 501826.9
$ws.Range("I5").Value = 768.1429000000001
$ws.Range("K5").Value = 2304.4287
$ws.Range("M5").Value = -2192.4287
$ws.Range("H33").Value = 47.857143
$ws.Range("J33").Value = 45
$ws.Range("L33").Value = 270
$ws.Range("N33").Value = -836
$ws.Range("H37").Value = 94677.62
$ws.Range("J37").Value = 94677.62
$ws.Range("L37").Value = 284032.86
$ws.Range("N37").Value = -284256.86
$ws.Range("H40").Value = 147.75
$ws.Range("I40").Value = 191
$ws.Range("J40").Value = 133.33333
$ws.Range("K40").Value = 764
$ws.Range("L40").Value = 533.33332
$ws.Range("M40").Value = -695
$ws.Range("N40").Value = -671.33332
$ws.Range("H113").Value = 1646.1111
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1830.7142
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 5492.142599999999
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -9832.142599999999
$ws.Range("H122").Value = 4522.952
$ws.Range("J122").Value = 4824.579
$ws.Range("L122").Value = 43421.211
$ws.Range("N122").Value = -48321.211
$ws.Range("H135").Value = 501826.9
$ws.Range("I135").Value = 768.1429000000001
$ws.Range("K135").Value = 6913.2861
$ws.Range("M135").Value = -4378.2861

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 13695.5
$ws.Range("I102").Value = 14578.4
$ws.Range("J102").Value = 6338
$ws.Range("K102").Value = 14578.4
$ws.Range("L102").Value = 6338
$ws.Range("M102").Value = -12956.4
$ws.Range("N102").Value = -9582
$ws.Range("H132").Value = 3809.5854
$ws.Range("I132").Value = 3839.3713
$ws.Range("K132").Value = 11518.1139
$ws.Range("M132").Value = -8988.1139
$ws.Range("H134").Value = 100326
$ws.Range("J134").Value = 100326
$ws.Range("L134").Value = 300978
$ws.Range("N134").Value = -306048
$ws.Range("H135").Value = 64093.43
$ws.Range("J135").Value = 64093.43
$ws.Range("L135").Value = 64093.43
$ws.Range("N135").Value = -74233.42999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 29790.234
$ws.Range("I7").Value = 45043.6
$ws.Range("J7").Value = 7999.7144
$ws.Range("K7").Value = 45043.6
$ws.Range("L7").Value = 7999.7144
$ws.Range("M7").Value = -44931.6
$ws.Range("N7").Value = -8223.714400000001
$ws.Range("H61").Value = 6187.143
$ws.Range("I61").Value = 2116.375
$ws.Range("K61").Value = 2116.375
$ws.Range("M61").Value = -1914.375
$ws.Range("H113").Value = 6187.143
$ws.Range("I113").Value = 2116.375
$ws.Range("K113").Value = 2116.375
$ws.Range("M113").Value = 53.625
$ws.Range("H122").Value = 7634.8887
$ws.Range("I122").Value = 8554.75
$ws.Range("K122").Value = 25664.25
$ws.Range("M122").Value = -23214.25
$ws.Range("H126").Value = 29790.234
$ws.Range("I126").Value = 45043.6
$ws.Range("J126").Value = 7999.7144
$ws.Range("K126").Value = 135130.8
$ws.Range("L126").Value = 23999.1432
$ws.Range("M126").Value = -132660.8
$ws.Range("N126").Value = -28939.1432
$ws.Range("H132").Value = 342223.5
$ws.Range("I132").Value = 535098.3
$ws.Range("K132").Value = 1605294.9
$ws.Range("M132").Value = -1602764.9

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 290462.66
$ws.Range("J62").Value = 10266.333
$ws.Range("L62").Value = 10266.333
$ws.Range("N62").Value = -11514.333
$ws.Range("H65").Value = 290462.66
$ws.Range("J65").Value = 10266.333
$ws.Range("L65").Value = 51331.665
$ws.Range("N65").Value = -57571.665
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H113").Value = 1761.88
$ws.Range("I113").Value = 988.8333
$ws.Range("J113").Value = 3749.7144
$ws.Range("K113").Value = 2966.4999
$ws.Range("L113").Value = 11249.1432
$ws.Range("M113").Value = -796.4998999999998
$ws.Range("N113").Value = -15589.1432
$ws.Range("H122").Value = 3980.7678
$ws.Range("I122").Value = 2538.7715
$ws.Range("K122").Value = 7616.314499999999
$ws.Range("M122").Value = -5166.314499999999
$ws.Range("H126").Value = 14139.648
$ws.Range("I126").Value = 17525.857
$ws.Range("K126").Value = 52577.571
$ws.Range("M126").Value = -50107.571
$ws.Range("H136").Value = 522938.1
$ws.Range("I136").Value = 556585.4399999999
$ws.Range("K136").Value = 1669756.32
$ws.Range("M136").Value = -1667206.32

